$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 13343792
$ws.Range("J88").Value = 11740.583
$ws.Range("L88").Value = 11740.583
$ws.Range("N88").Value = -12552.583
$ws.Range("H91").Value = 13343792
$ws.Range("J91").Value = 11740.583
$ws.Range("L91").Value = 11740.583
$ws.Range("N91").Value = -14548.583
$ws.Range("H106").Value = 18520880
$ws.Range("I106").Value = 23810848
$ws.Range("J106").Value = 5987.5
$ws.Range("K106").Value = 23810848
$ws.Range("L106").Value = 5987.5
$ws.Range("M106").Value = -23810217
$ws.Range("N106").Value = -7249.5
$ws.Range("H140").Value = 73461.53999999999
$ws.Range("I140").Value = 73333.336
$ws.Range("J140").Value = 75000
$ws.Range("K140").Value = 73333.336
$ws.Range("L140").Value = 75000
$ws.Range("M140").Value = -68153.336
$ws.Range("N140").Value = -85360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2004
$ws.Range("I45").Value = 2313.3333
$ws.Range("K45").Value = 2313.3333
$ws.Range("M45").Value = -1936.3333
$ws.Range("H61").Value = 13417.704
$ws.Range("I61").Value = 12753.207
$ws.Range("K61").Value = 12753.207
$ws.Range("M61").Value = -12541.207
$ws.Range("H88").Value = 46376.727
$ws.Range("I88").Value = 889
$ws.Range("J88").Value = 84283.164
$ws.Range("K88").Value = 889
$ws.Range("L88").Value = 84283.164
$ws.Range("M88").Value = -483
$ws.Range("N88").Value = -85095.164
$ws.Range("H91").Value = 46376.727
$ws.Range("I91").Value = 889
$ws.Range("J91").Value = 84283.164
$ws.Range("K91").Value = 889
$ws.Range("L91").Value = 84283.164
$ws.Range("M91").Value = 515
$ws.Range("N91").Value = -87091.164
$ws.Range("H122").Value = 4332.0415
$ws.Range("I122").Value = 2460
$ws.Range("J122").Value = 8878.429
$ws.Range("K122").Value = 7380
$ws.Range("L122").Value = 26635.287
$ws.Range("M122").Value = -4930
$ws.Range("N122").Value = -31535.287
$ws.Range("H136").Value = 13417.704
$ws.Range("I136").Value = 12753.207
$ws.Range("K136").Value = 38259.621
$ws.Range("M136").Value = -35709.621
$ws.Range("H140").Value = 110218.6
$ws.Range("J140").Value = 110218.6
$ws.Range("L140").Value = 110218.6
$ws.Range("N140").Value = -120578.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null
$ws.Range("H134").Value = 2537.0833
$ws.Range("I134").Value = 2087.1667
$ws.Range("J134").Value = 4786.6665
$ws.Range("K134").Value = 6261.500100000001
$ws.Range("L134").Value = 14359.9995
$ws.Range("M134").Value = -3726.500100000001
$ws.Range("N134").Value = -19429.9995
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1768.5454
$ws.Range("I16").Value = 1739.8889
$ws.Range("J16").Value = 1897.5
$ws.Range("K16").Value = 1739.8889
$ws.Range("L16").Value = 1897.5
$ws.Range("M16").Value = -1452.8889
$ws.Range("N16").Value = -2471.5
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H31").Value = 3963.8276
$ws.Range("J31").Value = 5097.7144
$ws.Range("L31").Value = 5097.7144
$ws.Range("N31").Value = -5687.7144
$ws.Range("H34").Value = 3963.8276
$ws.Range("J34").Value = 5097.7144
$ws.Range("L34").Value = 5097.7144
$ws.Range("N34").Value = -5501.7144
$ws.Range("H105").Value = 2273862.8
$ws.Range("I105").Value = 3247489.8
$ws.Range("J105").Value = 2066.6667
$ws.Range("K105").Value = 3247489.8
$ws.Range("L105").Value = 2066.6667
$ws.Range("M105").Value = -3245742.8
$ws.Range("N105").Value = -5560.6667
$ws.Range("H109").Value = 61590.168
$ws.Range("J109").Value = 61590.168
$ws.Range("L109").Value = 61590.168
$ws.Range("N109").Value = -63670.168
$ws.Range("H113").Value = 1768.5454
$ws.Range("I113").Value = 1739.8889
$ws.Range("J113").Value = 1897.5
$ws.Range("K113").Value = 1739.8889
$ws.Range("L113").Value = 1897.5
$ws.Range("M113").Value = 430.1111000000001
$ws.Range("N113").Value = -6237.5
$ws.Range("H122").Value = 5041
$ws.Range("I122").Value = 3235.3333
$ws.Range("J122").Value = 6395.25
$ws.Range("K122").Value = 9705.999899999999
$ws.Range("L122").Value = 19185.75
$ws.Range("M122").Value = -7255.999899999999
$ws.Range("N122").Value = -24085.75
$ws.Range("H132").Value = 9777.223
$ws.Range("I132").Value = 3916
$ws.Range("K132").Value = 11748
$ws.Range("M132").Value = -9218
$ws.Range("H134").Value = 2950.8
$ws.Range("I134").Value = 2842.9473
$ws.Range("K134").Value = 8528.841899999999
$ws.Range("M134").Value = -5993.841899999999
$ws.Range("H141").Value = 94107.48
$ws.Range("I141").Value = 42640
$ws.Range("J141").Value = 100162.47
$ws.Range("K141").Value = 42640
$ws.Range("L141").Value = 100162.47
$ws.Range("M141").Value = -37460
$ws.Range("N141").Value = -110522.47
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7382.6895
$ws.Range("I102").Value = 7096.1304
$ws.Range("K102").Value = 7096.1304
$ws.Range("M102").Value = -5474.1304
$ws.Range("H132").Value = 2428.842
$ws.Range("I132").Value = 2593.625
$ws.Range("K132").Value = 7780.875
$ws.Range("M132").Value = -5250.875
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4982.184
$ws.Range("I7").Value = 3834.6428
$ws.Range("K7").Value = 3834.6428
$ws.Range("M7").Value = -3722.6428
$ws.Range("H16").Value = 1918.375
$ws.Range("J16").Value = 5001
$ws.Range("L16").Value = 5001
$ws.Range("N16").Value = -5341
$ws.Range("H22").Value = 915.63635
$ws.Range("I22").Value = 591.3333
$ws.Range("K22").Value = 591.3333
$ws.Range("M22").Value = -296.3333
$ws.Range("H27").Value = 915.63635
$ws.Range("I27").Value = 591.3333
$ws.Range("K27").Value = 591.3333
$ws.Range("M27").Value = -484.3333
$ws.Range("H40").Value = 4068.6316
$ws.Range("I40").Value = 2881.8333
$ws.Range("K40").Value = 2881.8333
$ws.Range("M40").Value = -2745.8333
$ws.Range("H126").Value = 4982.184
$ws.Range("I126").Value = 3834.6428
$ws.Range("K126").Value = 11503.9284
$ws.Range("M126").Value = -9033.928400000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2536.1304
$ws.Range("I122").Value = 2321.55
$ws.Range("K122").Value = 6964.650000000001
$ws.Range("M122").Value = -4514.650000000001
